$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new activity description used by the new row entry.
$newActivity = "Design - Finished revising ER Diagrams, created complete ER Diagram"

# Fill in row 11 with the new progress-report entry (date, hours, activity).
# Match the date-formatted style already used by rows 3-10 in column A
# by copying A10's format down into A11 (keeps font/border, switches number format to date).
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A11").Value = 42791
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = $newActivity

# Update the active selection to reflect where the user left off (A12).
$ws.Range("A12").Select() | Out-Null
